# Add a new "Global ID" column (M) with a GUID-like value per data row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell - same bold/centered style (s="1") as the rest of row 1.
$ws.Range("M1").Value = "Global ID"
$ws.Range("M1").Font.Bold = $true
$ws.Range("M1").HorizontalAlignment = -4108
$ws.Range("M1").VerticalAlignment = -4108

# Data rows - entered in the same order the author typed them so that the
# shared-string table is appended in the same sequence.
$ws.Range("M2").Value = "{12BCC37E-15D7-4AA5-9287-7FBEE2C31483}"
$ws.Range("M4").Value = "{12BCC37E-15D7-IJNLKA-9287-ASDAS080122}"
$ws.Range("M3").Value = "{12BCC37E-15D7-ASDS-9287-POP21389123}"
$ws.Range("M5").Value = "{12BCC37E-15D7-4AA5-9287-ASDASD1389183}"
$ws.Range("M6").Value = "{12BCC37E-15D7-4AA5-9287-ASKPQ09121133}"

# Widen the new column to fit its content.
$ws.Columns.Item(13).AutoFit()

# Move / collapse the active selection onto the last filled cell.
[void]$ws.Range("M6").Select()

# Match the printed page setup recorded for this sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
